$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Remove the extra (empty) "Sheet2" / "Sheet3" worksheets.
# ---------------------------------------------------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets("Sheet2").Delete()
$wb.Worksheets("Sheet3").Delete()

# ---------------------------------------------------------------------------
# 2. "Sheet1" data sheet updates.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets("Sheet1")

# -- New row 38 ("basic RNN 1 layer") right after the existing row 37,
#    copying the row style used elsewhere in the workbook for this kind of
#    standalone label row (style index 40, taken from "official history"!A29).
$officialHistory = $wb.Worksheets("official history")
$officialHistory.Range("A29").Copy()
$ws1.Range("A38").PasteSpecial(-4122)
$ws1.Range("A38").Value = "basic RNN 1 layer"

$row38Vals = @(95.73,95.87,96.29,96.49,96.58,97.07,96.95,96.82,96.77,97.04,96.78,96.67,96.75)
$c = 2
foreach ($v in $row38Vals) {
    $ws1.Cells.Item(38, $c).Value = $v
    $c++
}

# -- Fill in previously-empty data cells in the existing LSTM/GRU/biGRU
#    comparison blocks.
$ws1.Cells.Item(43, 9).Value  = 73            # I43
$ws1.Cells.Item(43, 10).Value = 71.45         # J43
$ws1.Cells.Item(43, 11).Value = 71.46         # K43
$ws1.Cells.Item(43, 12).Value = 69.09         # L43
$ws1.Cells.Item(44, 9).Value  = 73.04         # I44
$ws1.Cells.Item(45, 15).Value = 73.81         # O45
$ws1.Cells.Item(46, 8).Value  = 72.38         # H46
$ws1.Cells.Item(51, 12).Value = 91.58         # L51
$ws1.Cells.Item(53, 4).Value  = 92.25         # D53
$ws1.Cells.Item(53, 5).Value  = 91.7          # E53
$ws1.Cells.Item(53, 6).Value  = 91.74         # F53
$ws1.Cells.Item(54, 6).Value  = 92.28         # F54
$ws1.Cells.Item(54, 7).Value  = 92.28         # G54
$ws1.Cells.Item(55, 4).Value  = 91.95         # D55

# -- New block: rows 58-63 ("s12" header + basic RNN/LSTM/GRU/biGRU rows),
#    mirroring the existing "m50" block's layout/format (rows 50-55).
$ws1.Range("A50:N50").Copy()
$ws1.Range("A58:N58").PasteSpecial(-4122)

$ws1.Range("B43:N46").Copy()
$ws1.Range("B59:N62").PasteSpecial(-4122)
$ws1.Range("B46:N46").Copy()
$ws1.Range("B63:N63").PasteSpecial(-4122)

$ws1.Range("A43:A46").Copy()
$ws1.Range("A59:A62").PasteSpecial(-4122)
$ws1.Range("A46").Copy()
$ws1.Range("A63").PasteSpecial(-4122)

$ws1.Range("A58").Value = "s12"
$ws1.Range("A59").Value = "basic RNN"
$ws1.Range("A60").Value = "LSTM"
$ws1.Range("A61").Value = "GRU"
$ws1.Range("A62").Value = "biGRU "
$ws1.Range("A63").Value = "biGRU + attention"

$row58Vals = @(8,24,40,56,72,88,104,120,128,144,160,176,256)
$c = 2
foreach ($v in $row58Vals) {
    $ws1.Cells.Item(58, $c).Value = $v
    $c++
}

# ---------------------------------------------------------------------------
# 3. View-state tweaks (scroll position / active selection) on a few sheets.
# ---------------------------------------------------------------------------
$wsHistNew = $wb.Worksheets("history_new")
$wsHistNew.Activate()
$wsHistNew.Range("N9:W9").Select()

$wsOfficial = $wb.Worksheets("official history")
$wsOfficial.Activate()
$wsOfficial.Range("B20:K20").Select()

$ws1.Activate()
$ws1.Range("A28").Select()
$ws1.Range("D40").Select()

Write-Output "edit complete"
